$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")
foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Insert a new row at position 14; existing rows 14-27 shift down to 15-28
    $ws.Rows.Item(14).Insert()

    # Fill in the newly inserted row 14 with the new event entry
    $ws.Cells.Item(14, 1).Value2 = 13
    $ws.Cells.Item(14, 1).Borders.LineStyle = 1
    $ws.Cells.Item(14, 2).NumberFormat = "@"
    $ws.Cells.Item(14, 2).Value2 = "2024-02-16"
    $ws.Cells.Item(14, 3).Value2 = "苏州·运动番only专区-Good jump ACG"
    $ws.Cells.Item(14, 4).Value2 = "金山南路288号 广电国际会展中心"
    $ws.Cells.Item(14, 5).Value2 = "2024.02.16 10:00-02.17 17:00"
    $ws.Cells.Item(14, 6).Value2 = 0
    $ws.Cells.Item(14, 7).Value2 = 25
    $ws.Cells.Item(14, 8).Value2 = "https://show.bilibili.com/platform/detail.html?id=81435"
    $ws.Cells.Item(14, 9).Value2 = "//i0.hdslb.com/bfs/openplatform/202401/gatL3YjP1706236832019.jpeg"

    # Renumber column A (row index counter) for all shifted rows 15-28 -> values 14-27
    for ($r = 15; $r -le 28; $r++) {
        $ws.Cells.Item($r, 1).Value2 = $r - 1
    }

    # Apply updated 'want-to-go' counter values (organic growth) in column F
    $ws.Cells.Item(3, 6).Value2 = 289
    $ws.Cells.Item(4, 6).Value2 = 11056
    $ws.Cells.Item(5, 6).Value2 = 10251
    $ws.Cells.Item(8, 6).Value2 = 724
    $ws.Cells.Item(9, 6).Value2 = 101
    $ws.Cells.Item(13, 6).Value2 = 9594
    $ws.Cells.Item(18, 6).Value2 = 12
    $ws.Cells.Item(19, 6).Value2 = 84
    $ws.Cells.Item(22, 6).Value2 = 10781
}

Write-Output "edit complete"
